# Insert a new row at position 134 (shifts existing rows 134-228 down to 135-229)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(134).Insert()

$row = 134
$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = 44762
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112032
$ws.Cells.Item($row, 7).Value = "Zapallo italiano"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 12000
$ws.Cells.Item($row, 12).Value = 13000
$ws.Cells.Item($row, 13).Value = 12500
$ws.Cells.Item($row, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 250
$ws.Cells.Item($row, 17).Value = 50
$ws.Cells.Item($row, 18).Value = "Hortaliza"

Write-Output "Row 134 inserted and populated"
